# Applies the diff: inserts two new price rows for Ciruela "Black Amber"
# (Especial and Primera, $/caja 15 kilos granel) at the top of the
# Macroferia Regional de Talca block, pushing the existing rows down by
# two positions (old row N -> new row N+2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62..85 down by two rows by inserting two new blank rows
# at position 62. Excel's Insert() copies formatting from the row
# above, matching the original D-column (date) style.
$ws.Rows.Item(62).Insert()
$ws.Rows.Item(62).Insert()

# New row 62: Ciruela, Black Amber, Especial
$ws.Range("A62").Value = 5
$ws.Range("B62").Value = "Macroferia Regional de Talca"
$ws.Range("C62").Value = "Maule"
$ws.Range("D62").Value = 44609
$ws.Range("E62").Value = 7
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100103
$ws.Range("H62").Value = "Frutos de hueso (carozo)"
$ws.Range("I62").Value = 100103002
$ws.Range("J62").Value = "Ciruela"
$ws.Range("K62").Value = "Black Amber"
$ws.Range("L62").Value = "Especial"
$ws.Range("M62").Value = 200
$ws.Range("N62").Value = 12000
$ws.Range("O62").Value = 12000
$ws.Range("P62").Value = 12000
$ws.Range("Q62").Value = "`$/caja 15 kilos granel"
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 800
$ws.Range("T62").Value = 15

# New row 63: Ciruela, Black Amber, Primera
$ws.Range("A63").Value = 5
$ws.Range("B63").Value = "Macroferia Regional de Talca"
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = 44609
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100103
$ws.Range("H63").Value = "Frutos de hueso (carozo)"
$ws.Range("I63").Value = 100103002
$ws.Range("J63").Value = "Ciruela"
$ws.Range("K63").Value = "Black Amber"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 200
$ws.Range("N63").Value = 10000
$ws.Range("O63").Value = 10000
$ws.Range("P63").Value = 10000
$ws.Range("Q63").Value = "`$/caja 15 kilos granel"
$ws.Range("R63").Value = "Región de O'Higgins"
$ws.Range("S63").Value = 667
$ws.Range("T63").Value = 15
